# Fixed update to excel issue
# 1. Rename the "Requested quantity" headers on the two existing sheets.
# 2. Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row.
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Forecast data rows.
$rows = @(
    @(44934.99999999999, 0,   -189.5107303698267, 151.8685076515984),
    @(44962.99999999999, 61,  -111.1244304905833, 249.4285208879284),
    @(44969.99999999999, 82,  -103.1752117619223, 274.1672709188894),
    @(44976.99999999999, 103, -87.6899618090756,  263.6955544330432),
    @(44983.99999999999, 124, -70.5559745447587,  297.9076461259442),
    @(44990.99999999999, 145, -25.57687131490281, 324.053235026794),
    @(44997.99999999999, 166, -14.14543585521915, 342.34012808962),
    @(45004.99999999999, 187, 20.35228997665739,  373.5167684266536),
    @(45011.99999999999, 208, 23.42695203672206,  383.7198508724479),
    @(45018.99999999999, 229, 47.46380378589409,  411.4557161211555),
    @(45025.99999999999, 250, 69.61609174111842,  427.9368207720275),
    @(45032.99999999999, 271, 91.55987673160946,  451.8785989245304),
    @(45039.99999999999, 292, 103.5175324082287,  470.4337582004144),
    @(45046.99999999999, 313, 134.0658270030061,  491.9911093704498)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match formatting of the other sheets: bold header row, date-formatted "ds" column.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A15").PasteSpecial(-4122)

$excel.CutCopyMode = 0
